{"js": "// Apply the NMCARS-PART-5233 edits:\n//  - Several paragraphs get pStyle \"List 2\" (replacing an explicit\n//    widowControl=0 pPr, or being added where there was no pPr at all).\n//  - One paragraph's pStyle changes from \"Normal w/line\" to \"List 1\".\n//  - Three paragraphs (\"(A) ...\", \"(B) ...\", \"(C) ...\") have their single\n//    run split in two: the \"(A)\"/\"(B)\"/\"(C)\" label becomes its own run,\n//    separate from the remaining sentence text.\n\nconst body = context.document.body;\n\nasync function findParagraph(needle) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const para = results.items[0].paragraphs.getFirst();\n  return para;\n}\n\nasync function setParagraphStyle(needle, styleName) {\n  const para = await findParagraph(needle);\n  para.style = styleName;\n  await context.sync();\n}\n\n// Split a paragraph's leading \"(A)\"-like label into its own run, separate\n// from the rest of the paragraph text, by re-inserting the whole paragraph\n// as OOXML with two <w:r> elements (one for the label, one for the rest).\nasync function splitLabelRun(needle, label, styleVal) {\n  const para = await findParagraph(needle);\n  para.load(\"text\");\n  await context.sync();\n\n  const fullText = para.text.replace(/\\r$/, \"\");\n  if (fullText.indexOf(label) !== 0) {\n    throw new Error(\"Label \" + JSON.stringify(label) + \" not at start of \" + JSON.stringify(fullText));\n  }\n  const rest = fullText.substring(label.length);\n\n  const pPr = styleVal ? ('<w:pPr><w:pStyle w:val=\"' + styleVal + '\"/></w:pPr>') : \"\";\n\n  // Word only emits xml:space=\"preserve\" when the run text actually has\n  // leading/trailing whitespace that needs preserving.\n  const needsPreserve = (s) => /^\\s|\\s$/.test(s);\n  const labelAttr = needsPreserve(label) ? ' xml:space=\"preserve\"' : \"\";\n  const restAttr = needsPreserve(rest) ? ' xml:space=\"preserve\"' : \"\";\n\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' + pPr +\n    '<w:r><w:t' + labelAttr + '>' + escapeXml(label) + '</w:t></w:r>' +\n    '<w:r><w:t' + restAttr + '>' + escapeXml(rest) + '</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  para.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// 1) widowControl=0 paragraphs -> pStyle List 2\nawait setParagraphStyle(\"When DON actions alleged\", \"List 2\");\nawait setParagraphStyle(\"In exceptional cases where disruption\", \"List 2\");\n\n// 2) (A)/(B)/(C) runs split, keeping their existing pStyle \"List 4\"\nawait splitLabelRun(\"the existence of a legal basis for entitlement\", \"(A)\", \"List4\");\nawait splitLabelRun(\"facts meeting the elements of proof required to support the basis\", \"(B)\", \"List4\");\nawait splitLabelRun(\"adequate factual support for the amounts claimed\", \"(C)\", \"List4\");\n\n// 3) paragraphs with no pPr at all -> add pStyle List 2\nawait setParagraphStyle(\"Analyze the applicability and adequacy of the contractor's legal theory\", \"List 2\");\nawait setParagraphStyle(\"Analyze and evaluate the presence and adequacy of evidentiary facts\", \"List 2\");\nawait setParagraphStyle(\"Analyze the applicability and adequacy of any affirmative defense\", \"List 2\");\nawait setParagraphStyle(\"Analyze and evaluate the presence of any counterclaims\", \"List 2\");\n\nawait setParagraphStyle(\"A legal determination that the contractor is entitled\", \"List 2\");\nawait setParagraphStyle(\"Sufficient technical, administrative, and audit analyses\", \"List 2\");\nawait setParagraphStyle(\"A determination by the contracting officer with respect to the amount\", \"List 2\");\n\nawait setParagraphStyle(\"the requirements of the claim;\", \"List 2\");\nawait setParagraphStyle(\"the projected date of settlement of the claim\", \"List 2\");\nawait setParagraphStyle(\"other pertinent information, including comments\", \"List 2\");\n\n// 4) pStyle Normalwline -> List1\nawait setParagraphStyle(\"General. The DON Office of the General Counsel\", \"List 1\");\n\nawait setParagraphStyle(\"A detailed narrative statement of facts\", \"List 2\");\nawait setParagraphStyle(\"An analysis and evaluation (classified as attorney\", \"List 2\");\nawait setParagraphStyle(\"The advisory report, if any, of the reviewing official\", \"List 2\");\n", "ps1": "# Apply the NMCARS-PART-5233 edits:\n#  - Several paragraphs get style \"List 2\" (replacing an explicit\n#    widowControl=0 pPr, or being added where there was no pPr at all).\n#  - One paragraph's style changes from \"Normal w/line\" to \"List 1\".\n#  - Three paragraphs (\"(A) ...\", \"(B) ...\", \"(C) ...\") have their single\n#    run split in two: the \"(A)\"/\"(B)\"/\"(C)\" label becomes its own run,\n#    separate from the remaining sentence text.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphByText($needle) {\n    $all = $d.Paragraphs\n    for ($i = 1; $i -le $all.Count; $i++) {\n        $p = $all.Item($i)\n        if ($p.Range.Text.Contains($needle)) {\n            return $p\n        }\n    }\n    throw \"Paragraph containing '$needle' not found\"\n}\n\nfunction Set-ParagraphStyle($needle, $styleName) {\n    $p = Find-ParagraphByText $needle\n    $p.Style = $styleName\n}\n\nfunction Split-LabelRun($needle, $label, $styleVal) {\n    $p = Find-ParagraphByText $needle\n    $paraRange = $p.Range\n\n    # Strip the trailing paragraph mark to get the plain text content.\n    $fullText = $paraRange.Text\n    if ($fullText.EndsWith([char]13)) {\n        $fullText = $fullText.Substring(0, $fullText.Length - 1)\n    }\n    if (-not $fullText.StartsWith($label)) {\n        throw \"Label '$label' not at start of '$fullText'\"\n    }\n    $rest = $fullText.Substring($label.Length)\n\n    $pPr = \"\"\n    if ($styleVal) {\n        $pPr = '<w:pPr><w:pStyle w:val=\"' + $styleVal + '\"/></w:pPr>'\n    }\n\n    # Word only emits xml:space=\"preserve\" when the run text actually has\n    # leading/trailing whitespace that needs preserving.\n    $labelAttr = \"\"\n    if ($label -match '^\\s' -or $label -match '\\s$') {\n        $labelAttr = ' xml:space=\"preserve\"'\n    }\n    $restAttr = \"\"\n    if ($rest -match '^\\s' -or $rest -match '\\s$') {\n        $restAttr = ' xml:space=\"preserve\"'\n    }\n\n    $labelEsc = $label.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n    $restEsc = $rest.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p>' + $pPr +\n        '<w:r><w:t' + $labelAttr + '>' + $labelEsc + '</w:t></w:r>' +\n        '<w:r><w:t' + $restAttr + '>' + $restEsc + '</w:t></w:r>' +\n        '</w:p>' +\n        '</w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $paraRange.Text = \"\"\n    $paraRange.InsertXML($xml) | Out-Null\n}\n\n# 1) widowControl=0 paragraphs -> style List 2\nSet-ParagraphStyle \"When DON actions alleged\" \"List 2\"\nSet-ParagraphStyle \"In exceptional cases where disruption\" \"List 2\"\n\n# 2) (A)/(B)/(C) runs split, keeping their existing style \"List 4\"\nSplit-LabelRun \"the existence of a legal basis for entitlement\" \"(A)\" \"List4\"\nSplit-LabelRun \"facts meeting the elements of proof required to support the basis\" \"(B)\" \"List4\"\nSplit-LabelRun \"adequate factual support for the amounts claimed\" \"(C)\" \"List4\"\n\n# 3) paragraphs with no pPr at all -> add style List 2\nSet-ParagraphStyle \"Analyze the applicability and adequacy of the contractor's legal theory\" \"List 2\"\nSet-ParagraphStyle \"Analyze and evaluate the presence and adequacy of evidentiary facts\" \"List 2\"\nSet-ParagraphStyle \"Analyze the applicability and adequacy of any affirmative defense\" \"List 2\"\nSet-ParagraphStyle \"Analyze and evaluate the presence of any counterclaims\" \"List 2\"\n\nSet-ParagraphStyle \"A legal determination that the contractor is entitled\" \"List 2\"\nSet-ParagraphStyle \"Sufficient technical, administrative, and audit analyses\" \"List 2\"\nSet-ParagraphStyle \"A determination by the contracting officer with respect to the amount\" \"List 2\"\n\nSet-ParagraphStyle \"the requirements of the claim;\" \"List 2\"\nSet-ParagraphStyle \"the projected date of settlement of the claim\" \"List 2\"\nSet-ParagraphStyle \"other pertinent information, including comments\" \"List 2\"\n\n# 4) style Normalwline -> List1\nSet-ParagraphStyle \"General. The DON Office of the General Counsel\" \"List 1\"\n\nSet-ParagraphStyle \"A detailed narrative statement of facts\" \"List 2\"\nSet-ParagraphStyle \"An analysis and evaluation (classified as attorney\" \"List 2\"\nSet-ParagraphStyle \"The advisory report, if any, of the reviewing official\" \"List 2\"\n"}
